$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price and Volume) to Text format so that
# numeric-looking strings (e.g. "1.00", "8.90") are preserved exactly
# as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "90.391.71"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "3.088.91"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "232.57"
$ws.Range("E5").Value = "  +6.82%  "
$ws.Range("D6").Value = "624.94"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "1.10"
$ws.Range("E7").Value = "  -4.24%  "
$ws.Range("D8").Value = "0.362"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.727"
$ws.Range("E10").Value = "  -4.84%  "
$ws.Range("B11").Value = "LidoStakedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D11").Value = "2.498.40"
$ws.Range("E11").Value = "  -20.92%  "
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "36.41"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "5.48"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").Value = "90.207.41"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "3.667.50"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "3.090.88"
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").Value = "3.77"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").Value = "0.0000211"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D22").Value = "439.23"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "5.55"
$ws.Range("E23").Value = "  +6.25%  "
$ws.Range("D24").Value = "8.90"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "7.57"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "5.69"
$ws.Range("E26").Value = "  -5.05%  "
$ws.Range("D27").Value = "89.02"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "12.27"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "9.45"
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "0.977"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("E34").Value = "  +15.63%  "
$ws.Range("D35").Value = "26.30"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.154"
$ws.Range("E36").Value = "  +6.67%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +2.92%  "
$ws.Range("D38").Value = "508.20"
$ws.Range("E38").Value = "  -3.90%  "
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "7.03"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "1.28"
$ws.Range("D42").Value = "0.0909"
$ws.Range("E42").Value = "  +4.18%  "
$ws.Range("E43").Value = "  +55.63%  "
$ws.Range("D44").Value = "0.411"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "1.90"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").Value = "150.82"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("D49").Value = "0.689"
$ws.Range("E49").Value = "  +5.34%  "
$ws.Range("D50").Value = "45.09"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").Value = "1.34"
$ws.Range("E51").Value = "  +0.05%  "
